# Apply cryptocurrency price/volume updates to Sheet1
# For each target cell: force a Text number format before assigning the
# value so that numeric-looking strings (e.g. "214.04") are stored as text
# (matching the original inlineStr/text cell type) instead of being
# auto-converted to a number by Excel. ClearFormats() afterwards restores
# the cell's default (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '89.300.01'
Set-TextValue "E2" '  +3.82%  '

Set-TextValue "D3" '3.289.56'
Set-TextValue "E3" '  -0.85%  '

Set-TextValue "E4" '  +0.01%  '

Set-TextValue "D5" '214.04'
Set-TextValue "E5" '  -2.50%  '

Set-TextValue "D6" '630.32'
Set-TextValue "E6" '  -0.69%  '

Set-TextValue "D7" '0.387'
Set-TextValue "E7" '  +18.76%  '

Set-TextValue "D8" '0.688'
Set-TextValue "E8" '  +15.70%  '

Set-TextValue "E9" '  -0.01%  '

Set-TextValue "D10" '3.284.42'
Set-TextValue "E10" '  -0.88%  '

Set-TextValue "D11" '0.578'
Set-TextValue "E11" '  -3.23%  '

Set-TextValue "E12" '  +12.25%  '

Set-TextValue "D13" '0.0000264'
Set-TextValue "E13" '  -4.56%  '

Set-TextValue "D14" '34.23'
Set-TextValue "E14" '  +0.44%  '

Set-TextValue "D15" '3.887.37'
Set-TextValue "E15" '  -0.97%  '

Set-TextValue "D16" '5.38'
Set-TextValue "E16" '  -0.14%  '

Set-TextValue "D17" '88.877.66'
Set-TextValue "E17" '  +3.97%  '

Set-TextValue "D18" '3.297.93'
Set-TextValue "E18" '  -0.33%  '

Set-TextValue "D19" '14.16'
Set-TextValue "E19" '  -3.06%  '

Set-TextValue "D20" '3.11'
Set-TextValue "E20" '  -1.57%  '

Set-TextValue "D21" '437.86'
Set-TextValue "E21" '  -1.18%  '

Set-TextValue "D22" '8.90'
Set-TextValue "E22" '  -2.39%  '

Set-TextValue "E23" '  +3.43%  '

Set-TextValue "D24" '7.40'
Set-TextValue "E24" '  +0.41%  '

Set-TextValue "D25" '12.37'
Set-TextValue "E25" '  +1.26%  '

Set-TextValue "D26" '5.27'
Set-TextValue "E26" '  -3.33%  '

Set-TextValue "D27" '3.449.99'
Set-TextValue "E27" '  -1.08%  '

Set-TextValue "D28" '77.03'
Set-TextValue "E28" '  -1.44%  '

Set-TextValue "D29" '0.0000134'
Set-TextValue "E29" '  +3.09%  '

Set-TextValue "E30" '  +0.00%  '

Set-TextValue "D31" '0.195'
Set-TextValue "E31" '  +15.44%  '

Set-TextValue "E32" '  +0.15%  '

Set-TextValue "D33" '577.54'
Set-TextValue "E33" '  -5.23%  '

Set-TextValue "D34" '8.89'
Set-TextValue "E34" '  -3.62%  '

Set-TextValue "E35" '  -9.18%  '

Set-TextValue "D36" '7.29'
Set-TextValue "E36" '  +13.26%  '

Set-TextValue "D37" '1.97'
Set-TextValue "E37" '  -3.36%  '

Set-TextValue "E38" '  -7.21%  '

Set-TextValue "D39" '22.68'
Set-TextValue "E39" '  -2.96%  '

Set-TextValue "D40" '21.84'
Set-TextValue "E40" '  +2.65%  '

Set-TextValue "D41" '0.999'
Set-TextValue "E41" '  -0.06%  '

Set-TextValue "E42" '  -3.86%  '

Set-TextValue "D43" '2.03'
Set-TextValue "E43" '  -0.81%  '

Set-TextValue "E44" '  -2.27%  '

Set-TextValue "E45" '  +0.09%  '

Set-TextValue "D46" '154.75'
Set-TextValue "E46" '  -2.36%  '

Set-TextValue "D47" '181.36'

Set-TextValue "D48" '45.06'
Set-TextValue "E48" '  -0.33%  '

Set-TextValue "E49" '  -3.60%  '

Set-TextValue "D50" '0.0692'
Set-TextValue "E50" '  +23.85%  '

Set-TextValue "B51" 'Stellar'
Set-TextValue "C51" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D51" '0.126'
Set-TextValue "E51" '  +0.31%  '
